$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply number formats to the new data columns (B:D numeric 0.00, E percent, F:G text)
$ws.Range("B2:D10").NumberFormat = "0.00"
$ws.Range("E2:E10").NumberFormat = "0%"
$ws.Range("F2:G10").NumberFormat = "@"

# Fill every data row (2-10) across columns B:G with the placeholder "N/A"
$dataRows = 2..10
foreach ($r in $dataRows) {
    $ws.Range("B$r").Value = "N/A"
    $ws.Range("C$r").Value = "N/A"
    $ws.Range("D$r").Value = "N/A"
    $ws.Range("E$r").Value = "N/A"
    $ws.Range("F$r").Value = "N/A"
    $ws.Range("G$r").Value = "N/A"
}

# Row 9 (Sharding) has a real TPS figure instead of "N/A"
$ws.Range("B9").Value = 13000

# Restore the cursor / selection to match the saved workbook state
[void]$ws.Range("C17").Select()
